# Update "想去人数" (want-to-go count) values in column F for both the
# "展览" and "全部类型" worksheets (they hold duplicate data).

$wb = $excel.ActiveWorkbook

$updates = @{
    "F4"  = 140
    "F5"  = 55
    "F6"  = 478
    "F7"  = 1354
    "F8"  = 495
    "F9"  = 96
    "F10" = 168
    "F12" = 169
    "F13" = 95
    "F14" = 142
    "F15" = 133
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
